$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# --- Update the F-column (time_taken) timestamps in "data" ---
$times = @(
    "2021-10-05 14:22:22.905038",
    "2021-10-05 14:22:22.905044",
    "2021-10-05 14:22:22.905046",
    "2021-10-05 14:22:22.905048",
    "2021-10-05 14:22:22.905050",
    "2021-10-05 14:22:22.905052",
    "2021-10-05 14:22:22.905054",
    "2021-10-05 14:22:22.905056",
    "2021-10-05 14:22:22.905058",
    "2021-10-05 14:22:22.905060",
    "2021-10-05 14:22:22.905062",
    "2021-10-05 14:22:22.905064",
    "2021-10-05 14:22:22.905066",
    "2021-10-05 14:22:22.905068",
    "2021-10-05 14:22:22.905070",
    "2021-10-05 14:22:22.905072",
    "2021-10-05 14:22:22.905074",
    "2021-10-05 14:22:22.905077",
    "2021-10-05 14:22:22.905079",
    "2021-10-05 14:22:22.905081",
    "2021-10-05 14:22:22.905083"
)
for ($i = 0; $i -lt $times.Length; $i++) {
    $row = $i + 2
    $dataSheet.Cells.Item($row, 6).Value = $times[$i]
}

# --- Add the new "metadata" sheet after "data" ---
$newSheet = $wb.Worksheets.Add($null, $dataSheet)
$newSheet.Name = "metadata"

# Reuse the bold/centered/bordered header style from data!B1:F1 (extend one
# extra column across to cover metadata's extra "panel_get_request" column)
$dataSheet.Range("B1:F1").Copy()
$newSheet.Range("B1:F1").PasteSpecial(-4122)
$dataSheet.Range("B1").Copy()
$newSheet.Range("G1").PasteSpecial(-4122)
# Reuse the "row-index" style (s=1) from data!A2 for metadata!A2
$dataSheet.Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)

# Header row
$newSheet.Range("B1").Value = "data_name"
$newSheet.Range("C1").Value = "data_id"
$newSheet.Range("D1").Value = "data_version"
$newSheet.Range("E1").Value = "data_version_created"
$newSheet.Range("F1").Value = "panel_query_time"
$newSheet.Range("G1").Value = "panel_get_request"

# Data row
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "Pulmonary arterial hypertension"
$newSheet.Range("C2").Value = 193
$newSheet.Range("D2").Value = "'2.16"
$newSheet.Range("D2").Style = "Normal"
$newSheet.Range("E2").Value = "2021-08-02T14:20:57.402885Z"
$newSheet.Range("F2").Value = "2021-10-05 14:22:22.902681"
$newSheet.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/193/?format=json"

# Keep "data" as the active sheet (matches original bookView activeTab=0)
$dataSheet.Activate()
$dataSheet.Range("A1").Select()

Write-Output "done"
